$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper characters that are awkward to embed literally / need building up
# ---------------------------------------------------------------------------
$mu  = [char]0x00B5     # µ (micro sign)
$pm  = [char]0x00B1     # ± (plus-minus)

# ---------------------------------------------------------------------------
# New row of data: second temperature sensor (TMP36)
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "SENSOR DE TEMPERATURA TMP36"
$ws.Range("C5").Value = "Medir la temperatura de los sistemas de control ambiental, protección térmica, control de procesos industriales, alarmas contra incendios, monitores de sistemas de potencia y gestión térmica de la CPU."
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = "2.7V - 5.5V"
$ws.Range("F5").Value = "-"
$ws.Range("G5").Value = "2°C"

$h5 = [string]::Concat("0.5 ", $mu, "A")
$ws.Range("H5").Value = $h5
$ws.Range("H5").Characters(5, 2).Font.Name = "Calibri"
$ws.Range("H5").Characters(5, 2).Font.Size = 11

$ws.Range("I5").Value = "-"
$ws.Range("J5").Value = "-"
$ws.Range("K5").Value = "Análogo "
$ws.Range("L5").Value = "10mV/°C"
$ws.Range("M5").Value = "desde -50°C y 125°C"
$ws.Range("N5").Value2 = 5900
$ws.Range("O5").Value = "https://www.amazon.com/-/es/TMP36-Sensor-de-temperatura/dp/B00JYQAIBM"

# Row height for the new row (tall, to fit the wrapped "Aplicación" text)
$ws.Rows.Item(5).RowHeight = 240

# ---------------------------------------------------------------------------
# Formatting for row 5
# ---------------------------------------------------------------------------
$row5 = $ws.Range("B5:O5")
$row5.Borders.LineStyle = 1
$row5.HorizontalAlignment = -4108
$row5.VerticalAlignment = -4108
$row5.WrapText = $true

# H5 (consumption) keeps the "no wrap" look used by the rest of that column
$ws.Range("H5").WrapText = $false

# N5 (price) uses the same currency-like number format as N4
$ws.Range("N5").NumberFormat = '_-"$"\ * #,##0.00_-;\-"$"\ * #,##0.00_-;_-"$"\ * "-"??_-;_-@_-'

# ---------------------------------------------------------------------------
# Row 6 (now the trailing blank row) switches from the "no-wrap" look to the
# "wrap" look used elsewhere in the table.
# ---------------------------------------------------------------------------
$row6 = $ws.Range("B6:O6")
$row6.Borders.LineStyle = 1
$row6.HorizontalAlignment = -4108
$row6.VerticalAlignment = -4108
$row6.WrapText = $true

# ---------------------------------------------------------------------------
# Second annotation textbox (±) next to the new "Precisión" value in G5
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.AddTextbox(1, 473.25, 279, 10.81, 13.56)
$shp.Name = "CuadroTexto 2"
$shp.TextFrame.Characters().Text = $pm

# ---------------------------------------------------------------------------
# View settings
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 78
$ws.Range("O5").Select()
